$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: new TIP contact status record.
$ws.Range("A6").Value = "DNET COMMUNICATIONS"

# Empty text cell (matches the blank BBM_STD cells elsewhere in the sheet).
$ws.Range("B6").Value = "'"
$ws.Range("B6").ClearFormats()

$ws.Range("C6").Value = "OS"

# Account number must stay text (not be auto-converted to a number).
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "9028213376"
$ws.Range("D6").ClearFormats()

$ws.Range("E6").Value = "2025-12-02 14:28"

# Empty text cell (matches the blank LAST_WHATSAPP_TIME cells elsewhere).
$ws.Range("F6").Value = "'"
$ws.Range("F6").ClearFormats()

$ws.Range("G6").Value = "2025-12"
